$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.179.48"
$ws.Range("E2").Value = "  +1.74%  "
$ws.Range("D3").Value = "1.654.88"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'219.50"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("E7").Value = "  -0.21%  "
$ws.Range("D8").Value = "'0.255"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "'0.0625"
$ws.Range("E9").Value = "  +0.54%  "
$ws.Range("D10").Value = "'19.93"
$ws.Range("E10").Value = "  +4.76%  "
$ws.Range("D11").Value = "'0.0849"
$ws.Range("E11").Value = "  +0.79%  "
$ws.Range("D12").Value = "1.891.47"
$ws.Range("E12").Value = "  +1.47%  "
$ws.Range("D13").Value = "1.666.65"
$ws.Range("E13").Value = "  +2.64%  "
$ws.Range("D14").Value = "'4.17"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'0.531"
$ws.Range("E15").Value = "  +1.31%  "
$ws.Range("D16").Value = "'67.14"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "27.192.78"
$ws.Range("E17").Value = "  +1.91%  "
$ws.Range("E18").Value = "  +1.09%  "
$ws.Range("D19").Value = "'221.81"
$ws.Range("E19").Value = "  +5.29%  "
$ws.Range("E20").Value = "  -0.36%  "
$ws.Range("D21").Value = "'4.44"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("D22").Value = "'6.71"
$ws.Range("E22").Value = "  +8.56%  "
$ws.Range("E23").Value = "  +3.18%  "
$ws.Range("D24").Value = "'9.26"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'146.96"
$ws.Range("E25").Value = "  +1.03%  "
$ws.Range("E26").Value = "  -0.31%  "
$ws.Range("D27").Value = "'7.44"
$ws.Range("E27").Value = "  +5.30%  "
$ws.Range("E28").Value = "  +1.54%  "
$ws.Range("D29").Value = "'16.02"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").Value = "'0.0514"
$ws.Range("E30").Value = "  +1.88%  "
$ws.Range("E31").Value = "  +1.05%  "
$ws.Range("D32").Value = "'3.42"
$ws.Range("E32").Value = "  +2.11%  "
$ws.Range("D33").Value = "'3.01"
$ws.Range("E33").Value = "  +1.31%  "
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").Value = "1.258.23"
$ws.Range("E35").Value = "  -1.43%  "
$ws.Range("E36").Value = "  +1.10%  "
$ws.Range("D37").Value = "'0.0178"
$ws.Range("E37").Value = "  +1.65%  "
$ws.Range("D38").Value = "'0.536"
$ws.Range("E38").Value = "  +1.08%  "
$ws.Range("D39").Value = "'0.837"
$ws.Range("E39").Value = "  +3.73%  "
$ws.Range("E40").Value = "  -0.29%  "
$ws.Range("D41").Value = "'0.816"
$ws.Range("E41").Value = "  +1.90%  "
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("D43").Value = "1.801.00"
$ws.Range("E43").Value = "  +1.46%  "
$ws.Range("E44").Value = "  -3.67%  "
$ws.Range("D45").Value = "'61.85"
$ws.Range("E45").Value = "  +2.18%  "
$ws.Range("D46").Value = "'91.81"
$ws.Range("E46").Value = "  +0.12%  "
$ws.Range("E47").Value = "  +2.42%  "
$ws.Range("D48").Value = "'0.0516"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("E49").Value = "  +2.47%  "
$ws.Range("E50").Value = "  +1.87%  "
$ws.Range("E51").Value = "  +0.44%  "
